$wb = $excel.ActiveWorkbook

# --- Rename Sheet1 -> fund_details, move selection off C6 to D8, drop tab-selected ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "fund_details"
$ws1.Range("D8").Select() | Out-Null

# --- Add new sheet "fund_map" right after fund_details ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "fund_map"

# Scratch cell used to materialise literal-leading-apostrophe strings without
# triggering Excel's quote-prefix (text-qualifier) autoformat.
$tmp = $ws2.Cells.Item(200,1)

# --- Populate the raw-to-normalized fund manager mapping table ---
$ws2.Cells.Item(1,1).Value = "fundManager_raw"
$ws2.Cells.Item(1,2).Value = "fundManager"
$tmp.Formula = "=CHAR(39)&""Azzad Funds""&CHAR(39)"
$tmp.Copy() | Out-Null
$ws2.Cells.Item(2,1).PasteSpecial(-4163) | Out-Null
$ws2.Cells.Item(2,2).Value = "Azzad"
$ws2.Cells.Item(3,1).Value = " 'BARONFUNDS'"
$ws2.Cells.Item(3,2).Value = "Baron"
$ws2.Cells.Item(4,1).Value = " 'BlackRock-Advised Funds'"
$ws2.Cells.Item(4,2).Value = "BlackRock"
$ws2.Cells.Item(5,1).Value = " 'BlackRock-advised Funds'"
$ws2.Cells.Item(5,2).Value = "BlackRock"
$ws2.Cells.Item(6,1).Value = " 'Brighthouse Funds Trust I'"
$ws2.Cells.Item(6,2).Value = "Brighthouse"
$ws2.Cells.Item(7,1).Value = " 'Brighthouse Funds Trust II'"
$ws2.Cells.Item(7,2).Value = "Brighthouse"
$ws2.Cells.Item(8,1).Value = " 'Brinker Capital Destinations Trust'"
$ws2.Cells.Item(8,2).Value = "Brinker Capital"
$ws2.Cells.Item(9,1).Value = " 'CRMCFNDGRP'"
$ws2.Cells.Item(9,2).Value = "Capital Group"
$ws2.Cells.Item(10,1).Value = " 'Clipper Funds Trust'"
$ws2.Cells.Item(10,2).Value = "Clipper"
$ws2.Cells.Item(11,1).Value = " 'Davis Funds'"
$ws2.Cells.Item(11,2).Value = "Davis "
$ws2.Cells.Item(12,1).Value = " 'Delaware Funds by Macquarie'"
$ws2.Cells.Item(12,2).Value = "Macquarie"
$ws2.Cells.Item(13,1).Value = " 'FIRSTPACAD'"
$ws2.Cells.Item(13,2).Value = "First Pacific Advisors"
$ws2.Cells.Item(14,1).Value = " 'Federated Hermes Funds'"
$ws2.Cells.Item(14,2).Value = "Federated Hermes"
$ws2.Cells.Item(15,1).Value = " 'Fidelity Group of Funds'"
$ws2.Cells.Item(15,2).Value = "Fidelity"
$ws2.Cells.Item(16,1).Value = " 'Franklin Templeton Group of Funds'"
$ws2.Cells.Item(16,2).Value = "Franklin Templeton"
$ws2.Cells.Item(17,1).Value = " 'GREAT-WEST FUNDS INC'"
$ws2.Cells.Item(17,2).Value = "Great-West"
$ws2.Cells.Item(18,1).Value = " 'Goldman Sachs Fund Complex'"
$ws2.Cells.Item(18,2).Value = "Goldman Sachs"
$ws2.Cells.Item(19,1).Value = " 'Guggenheim'"
$ws2.Cells.Item(19,2).Value = "Guggenheim"
$ws2.Cells.Item(20,1).Value = " 'HARTFORD FUNDS'"
$ws2.Cells.Item(20,2).Value = "Hartford Funds"
$ws2.Cells.Item(21,1).Value = " 'INVESCOFDS'"
$ws2.Cells.Item(21,2).Value = "Invesco"
$ws2.Cells.Item(22,1).Value = " 'JNL Fund Complex'"
$ws2.Cells.Item(22,2).Value = "Jackson National"
$ws2.Cells.Item(23,1).Value = " 'John Hancock Group of Funds'"
$ws2.Cells.Item(23,2).Value = "John Hancock"
$ws2.Cells.Item(24,1).Value = " 'LINCOLNTRS'"
$ws2.Cells.Item(24,2).Value = "Lincoln Investment"
$ws2.Cells.Item(25,1).Value = " 'MASSMUTUAL FUNDS'"
$ws2.Cells.Item(25,2).Value = "MassMutual"
$ws2.Cells.Item(26,1).Value = " 'Neuberger Berman'"
$ws2.Cells.Item(26,2).Value = "Neuberger Berman"
$ws2.Cells.Item(27,1).Value = " 'PENN SERIES FUNDS INC'"
$ws2.Cells.Item(27,2).Value = "PennMutual"
$ws2.Cells.Item(28,1).Value = " 'Principal Funds'"
$ws2.Cells.Item(28,2).Value = "Principal Funds"
$ws2.Cells.Item(29,1).Value = " 'Publicly registered funds managed by Blackstone Alternative Asset Management L.P. Blackstone Alternative Investment Advisors LLC'"
$ws2.Cells.Item(29,2).Value = "Blackstone"
$ws2.Cells.Item(30,1).Value = " 'SUNAMERICA'"
$ws2.Cells.Item(30,2).Value = "Sun America"
$ws2.Cells.Item(31,1).Value = " 'Selected Funds'"
$ws2.Cells.Item(31,2).Value = "Selected Funds"
$ws2.Cells.Item(32,1).Value = " 'TROWEPRICE'"
$ws2.Cells.Item(32,2).Value = "T. Rowe Price"
$ws2.Cells.Item(33,1).Value = " 'VALIC Company'"
$ws2.Cells.Item(33,2).Value = "VALIC"
$ws2.Cells.Item(34,1).Value = " 'Voya mutual funds'"
$ws2.Cells.Item(34,2).Value = "Voya"

$tmp.ClearContents() | Out-Null

# --- Column widths (closest achievable approximation of the authored widths) ---
$ws2.Columns.Item(1).ColumnWidth = 31.917
$ws2.Columns.Item(2).ColumnWidth = 25.251

# --- Select/activate fund_map as the active sheet+cell (becomes tabSelected) ---
$ws2.Range("B10").Select() | Out-Null
$ws2.Activate() | Out-Null

Write-Host "edit complete"
